$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 807.7037
$ws.Range("I19").Value = 247.41176
$ws.Range("J19").Value = 1760.2
$ws.Range("K19").Value = 247.41176
$ws.Range("L19").Value = 1760.2
$ws.Range("M19").Value = -72.41175999999999
$ws.Range("N19").Value = -2110.2
$ws.Range("H43").Value = 9279.286
$ws.Range("I43").Value = 4540
$ws.Range("J43").Value = 11912.223
$ws.Range("K43").Value = 4540
$ws.Range("L43").Value = 11912.223
$ws.Range("M43").Value = -4471
$ws.Range("N43").Value = -12050.223
$ws.Range("H51").Value = 3065.5293
$ws.Range("I51").Value = 3500
$ws.Range("J51").Value = 3007.6
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 3007.6
$ws.Range("M51").Value = -3016
$ws.Range("N51").Value = -3975.6
$ws.Range("H137").Value = 8334777
$ws.Range("I137").Value = 1202.5385
$ws.Range("K137").Value = 3607.6155
$ws.Range("M137").Value = -1057.6155
$ws.Range("H138").Value = 6758200
$ws.Range("J138").Value = 16668447
$ws.Range("L138").Value = 50005341
$ws.Range("N138").Value = -50015621

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11906995
$ws.Range("I61").Value = 12197361
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 12197361
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -12197149
$ws.Range("N61").Value = -2424
$ws.Range("H74").Value = 13890917
$ws.Range("I74").Value = 17242702
$ws.Range("J74").Value = 4948.857
$ws.Range("K74").Value = 17242702
$ws.Range("L74").Value = 4948.857
$ws.Range("M74").Value = -17241828
$ws.Range("N74").Value = -6696.857
$ws.Range("H77").Value = 13890917
$ws.Range("I77").Value = 17242702
$ws.Range("J77").Value = 4948.857
$ws.Range("K77").Value = 86213510
$ws.Range("L77").Value = 24744.285
$ws.Range("M77").Value = -86209142
$ws.Range("N77").Value = -33480.285
$ws.Range("H132").Value = 19235152
$ws.Range("I132").Value = 50005404
$ws.Range("J132").Value = 3744
$ws.Range("K132").Value = 150016212
$ws.Range("L132").Value = 11232
$ws.Range("M132").Value = -150013682
$ws.Range("N132").Value = -16292
$ws.Range("H136").Value = 11906995
$ws.Range("I136").Value = 12197361
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 36592083
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -36589533
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39263.332
$ws.Range("J35").Value = 39263.332
$ws.Range("L35").Value = 39263.332
$ws.Range("N35").Value = -39883.332
$ws.Range("H107").Value = 3863.75
$ws.Range("I107").Value = 4411
$ws.Range("J107").Value = 2222
$ws.Range("K107").Value = 4411
$ws.Range("L107").Value = 2222
$ws.Range("M107").Value = -2491
$ws.Range("N107").Value = -6062
$ws.Range("H134").Value = 6043.6
$ws.Range("I134").Value = 3527.2
$ws.Range("J134").Value = 8560
$ws.Range("K134").Value = 10581.6
$ws.Range("L134").Value = 25680
$ws.Range("M134").Value = -8046.599999999999
$ws.Range("N134").Value = -30750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1450
$ws.Range("J16").Value = 1675
$ws.Range("L16").Value = 1675
$ws.Range("N16").Value = -2249
$ws.Range("H31").Value = 5750049
$ws.Range("I31").Value = 2950.5293
$ws.Range("J31").Value = 47621764
$ws.Range("K31").Value = 2950.5293
$ws.Range("L31").Value = 47621764
$ws.Range("M31").Value = -2655.5293
$ws.Range("N31").Value = -47622354
$ws.Range("H34").Value = 5750049
$ws.Range("I34").Value = 2950.5293
$ws.Range("J34").Value = 47621764
$ws.Range("K34").Value = 2950.5293
$ws.Range("L34").Value = 47621764
$ws.Range("M34").Value = -2748.5293
$ws.Range("N34").Value = -47622168
$ws.Range("H53").Value = 32592
$ws.Range("J53").Value = 32592
$ws.Range("L53").Value = 32592
$ws.Range("N53").Value = -33806
$ws.Range("H58").Value = 2630
$ws.Range("I58").Value = 1165.4615
$ws.Range("J58").Value = 5803.1665
$ws.Range("K58").Value = 1165.4615
$ws.Range("L58").Value = 5803.1665
$ws.Range("M58").Value = -962.4614999999999
$ws.Range("N58").Value = -6209.1665
$ws.Range("H113").Value = 1450
$ws.Range("J113").Value = 1675
$ws.Range("L113").Value = 1675
$ws.Range("N113").Value = -6015
$ws.Range("H134").Value = 555265.4
$ws.Range("I134").Value = 1330.5454
$ws.Range("K134").Value = 3991.6362
$ws.Range("M134").Value = -1456.6362
$ws.Range("H136").Value = 2630
$ws.Range("I136").Value = 1165.4615
$ws.Range("J136").Value = 5803.1665
$ws.Range("K136").Value = 3496.3845
$ws.Range("L136").Value = 17409.4995
$ws.Range("M136").Value = -946.3844999999997
$ws.Range("N136").Value = -22509.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 600.5263
$ws.Range("I5").Value = 365.45456
$ws.Range("J5").Value = 923.75
$ws.Range("K5").Value = 1096.36368
$ws.Range("L5").Value = 2771.25
$ws.Range("M5").Value = -984.3636799999999
$ws.Range("N5").Value = -2995.25
$ws.Range("H68").Value = 1113
$ws.Range("I68").Value = 567
$ws.Range("J68").Value = 1440.6
$ws.Range("K68").Value = 1701
$ws.Range("L68").Value = 4321.799999999999
$ws.Range("M68").Value = -890
$ws.Range("N68").Value = -5943.799999999999
$ws.Range("H71").Value = 1113
$ws.Range("I71").Value = 567
$ws.Range("J71").Value = 1440.6
$ws.Range("K71").Value = 5103
$ws.Range("L71").Value = 12965.4
$ws.Range("M71").Value = -1047
$ws.Range("N71").Value = -21077.4
$ws.Range("H113").Value = 834.9655
$ws.Range("I113").Value = 512.9091
$ws.Range("K113").Value = 1538.7273
$ws.Range("M113").Value = 631.2727
$ws.Range("H132").Value = 663.4286
$ws.Range("I132").Value = 447.16666
$ws.Range("J132").Value = 825.625
$ws.Range("K132").Value = 4024.49994
$ws.Range("L132").Value = 7430.625
$ws.Range("M132").Value = -1494.49994
$ws.Range("N132").Value = -12490.625
$ws.Range("H135").Value = 600.5263
$ws.Range("I135").Value = 365.45456
$ws.Range("J135").Value = 923.75
$ws.Range("K135").Value = 3289.09104
$ws.Range("L135").Value = 8313.75
$ws.Range("M135").Value = -754.0910400000002
$ws.Range("N135").Value = -13383.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 6213.5
$ws.Range("I107").Value = 20000
$ws.Range("J107").Value = 1618
$ws.Range("K107").Value = 20000
$ws.Range("L107").Value = 1618
$ws.Range("M107").Value = -18080
$ws.Range("N107").Value = -5458
$ws.Range("H113").Value = 251481
$ws.Range("J113").Value = 1956.5
$ws.Range("L113").Value = 1956.5
$ws.Range("N113").Value = -6296.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6005.5557
$ws.Range("I7").Value = 6054.5454
$ws.Range("J7").Value = 5928.5713
$ws.Range("K7").Value = 6054.5454
$ws.Range("L7").Value = 5928.5713
$ws.Range("M7").Value = -5942.5454
$ws.Range("N7").Value = -6152.5713
$ws.Range("H126").Value = 6005.5557
$ws.Range("I126").Value = 6054.5454
$ws.Range("J126").Value = 5928.5713
$ws.Range("K126").Value = 18163.6362
$ws.Range("L126").Value = 17785.7139
$ws.Range("M126").Value = -15693.6362
$ws.Range("N126").Value = -22725.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1790.5333
$ws.Range("I113").Value = 428.66666
$ws.Range("J113").Value = 3833.3333
$ws.Range("K113").Value = 1285.99998
$ws.Range("L113").Value = 11499.9999
$ws.Range("M113").Value = 884.0000199999999
$ws.Range("N113").Value = -15839.9999
$ws.Range("H132").Value = 1268.3572
$ws.Range("I132").Value = 914.7646999999999
$ws.Range("J132").Value = 1814.8182
$ws.Range("K132").Value = 2744.2941
$ws.Range("L132").Value = 5444.4546
$ws.Range("M132").Value = -214.2941000000001
$ws.Range("N132").Value = -10504.4546
